# daily auto push: 2026-02-28 02:38 UTC
# Insert a new data row for 2026/02/28 (Sat), hour 7, ranking 27 at row 875,
# pushing the existing 2026/12/29..2027/01/05 rows down by one (through row 917).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 875 (the rest of the sheet
# below shifts down automatically, same as Excel's Rows.Insert()).
$ws.Rows.Item(875).Insert()

# Populate the new row. The date column stores plain text (e.g. "2026/12/29")
# in this workbook, not a real date value, so prefix with a leading
# apostrophe to force text entry and avoid Excel's automatic date parsing;
# then reset the style back to "Normal" so the cell doesn't pick up the
# quote-prefix style and matches the rest of the sheet.
$ws.Range("A875").Value = "'2026/02/28"
$ws.Range("A875").Style = "Normal"
$ws.Range("B875").Value = "土"
$ws.Range("C875").Value = 7
$ws.Range("D875").Value = 27
